$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Shrub" to "CSS" in the vegetation treatment labels.
# Resulting cell values (by position, matching the target sharedStrings order):
#   A1: groups                (unchanged)
#   A2: CSS x Reduced         (was "Grassland x Reduced")
#   A3: Grassland x Ambient   (was "Shrub x Ambient")
#   A4: CSS x Ambient         (was "Grassland x Ambient")
#   A5: Grassland x Reduced   (was "Shrub x Reduced")

$ws.Range("A2").Value = "CSS x Reduced"
$ws.Range("A3").Value = "Grassland x Ambient"
$ws.Range("A4").Value = "CSS x Ambient"
$ws.Range("A5").Value = "Grassland x Reduced"
